# 自动更新Excel文件 - 2025-12-01 23:12:39
# For each data row: decrement remaining days (column E) by 1.
# If remaining days was already down to 1 (last day), the item is treated as
# replenished: remaining resets to the total (column D) and the start date
# (column F) rolls to the next day. Rows whose start date isn't a valid
# yyyymmdd date are left untouched (mirrors the source data's malformed date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$todayYear  = 2025
$todayMonth = 12
$todayDay   = 1

function Test-ValidYmd($value) {
    $s = [string]([int64]$value)
    if ($s.Length -ne 8) { return $false }
    $y = [int]$s.Substring(0,4)
    $m = [int]$s.Substring(4,2)
    $d = [int]$s.Substring(6,2)
    if ($m -lt 1 -or $m -gt 12) { return $false }
    if ($d -lt 1 -or $d -gt 31) { return $false }
    try {
        $dt = Get-Date -Year $y -Month $m -Day $d
        return $true
    } catch {
        return $false
    }
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $totalCell     = $ws.Cells.Item($r, 4)
    $remainingCell = $ws.Cells.Item($r, 5)
    $startCell     = $ws.Cells.Item($r, 6)

    $total     = $totalCell.Value2
    $remaining = $remainingCell.Value2
    $start     = $startCell.Value2

    if ($null -eq $total -or $null -eq $remaining -or $null -eq $start) {
        continue
    }

    if (-not (Test-ValidYmd $start)) {
        continue
    }

    if ([int]$remaining -eq 1) {
        $remainingCell.Value2 = [int]$total
        $nextStart = (Get-Date -Year $todayYear -Month $todayMonth -Day $todayDay).AddDays(1)
        $startCell.Value2 = [int]$nextStart.ToString("yyyyMMdd")
    } else {
        $remainingCell.Value2 = [int]$remaining - 1
    }
}
